$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-11-21 Thursday" "2024-11-22 Friday"

Replace-Text "530×8=4240" "921×4=3684"
Replace-Text "743×9=6687" "289×3=867"
Replace-Text "797×2=1594" "148×9=1332"
Replace-Text "877×8=7016" "152×2=304"
Replace-Text "209×3=627" "956×2=1912"

Replace-Text "229×9=2061" "691×7=4837"
Replace-Text "382×9=3438" "601×4=2404"
Replace-Text "902×6=5412" "267×4=1068"
Replace-Text "339×3=1017" "665×3=1995"
Replace-Text "597×5=2985" "677×9=6093"

Replace-Text "145×2=290" "647×7=4529"
Replace-Text "528×8=4224" "990×4=3960"
Replace-Text "139×2=278" "221×9=1989"
Replace-Text "549×6=3294" "539×2=1078"
Replace-Text "161×7=1127" "709×8=5672"

Replace-Text "551×3=1653" "383×9=3447"
Replace-Text "140×3=420" "638×3=1914"
Replace-Text "657×8=5256" "754×6=4524"
Replace-Text "811×4=3244" "132×8=1056"
Replace-Text "163×5=815" "922×5=4610"

Replace-Text "855×6=5130" "221×2=442"
Replace-Text "154×3=462" "741×5=3705"
Replace-Text "881×6=5286" "137×3=411"
Replace-Text "181×9=1629" "355×6=2130"
Replace-Text "309×2=618" "811×6=4866"
